# "Generate Report for Archive"
# - Update the localization status text from "Ready for handoff" to "In Translation"
#   (this text is a shared string referenced from the Overview sheet's zh-cn/de-de
#   status columns as well as the Status column on each per-language sheet).
# - Shrink the "Status" column width on the Overview sheet (columns E & F, i.e. the
#   zh-cn / de-de columns) and on the zh-cn / de-de sheets (column C, the Status
#   column) from ~17.22 to ~13.41 characters.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The COM layer stores column widths as "characters" but rounds to whole pixels
# (Maximum Digit Width = 6 for this runtime) with a fixed 5-pixel padding, i.e.
# storedWidth = round((ColumnWidth + 5/6) * 6) / 6. Back out the ColumnWidth to
# feed in order to land as close as possible to the target stored width.
$targetStoredWidth = 13.4101845877511
$columnWidthToSet = $targetStoredWidth - (5.0 / 6.0)

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # NOTE: compare with the literal on the left-hand side. Value2 can come
            # back as a non-string (e.g. a boolean for cells literally containing
            # "True"/"False"), and PowerShell's -eq coerces the right-hand side to
            # the left-hand side's type, which would otherwise give false matches.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }

    if ($ws.Name -eq "Overview") {
        $ws.Range("E1").ColumnWidth = $columnWidthToSet
        $ws.Range("F1").ColumnWidth = $columnWidthToSet
    } elseif ($ws.Name -eq "zh-cn" -or $ws.Name -eq "de-de") {
        $ws.Range("C1").ColumnWidth = $columnWidthToSet
    }
}
